$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1: update meanEMG-related header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: legmaxROM value moved from D2 to C2
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 30.373805491377226

# Row 3: remove stray values from B3 and C3
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Update selection to reflect the new active range
$ws.Range("B1:E3").Select()
